$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for handoff"
#
# The 937e6c38-... item is now ready to be handed off again, so it moves back
# to the top of the per-language sheets (swapping places with the
# f8972293-... item) and its status / handoff timestamp are refreshed.
# ---------------------------------------------------------------------------

# ---- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.md"
$wsOverview.Range("A3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/e2e/937e6c38-aff9-4180-b7f0-59051d5eb1b8.md", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/e2e/f8972293-0ef2-427e-b81c-380ee44cbdd1.md", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/.localization-config", "", "", ".localization-config") | Out-Null

# ---- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.md"
$wsZh.Range("C2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.zh-cn.xlf"
$wsZh.Range("E2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.md"
$wsZh.Range("F2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.zh-cn.xlf"

$wsZh.Range("A3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md"
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-01-26 11:53:19"
$wsZh.Range("E3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md"
$wsZh.Range("F3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/e2e/937e6c38-aff9-4180-b7f0-59051d5eb1b8.md", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/417d83acef18db16e4504269791455a1b7e15ec5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.zh-cn.xlf", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c9255525979e9db716f4edd38e1c5885a9e23a62/e2e/937e6c38-aff9-4180-b7f0-59051d5eb1b8.md", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74025245821fd71bdd61467e2aa2585ffc1af607/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.zh-cn.xlf", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/e2e/f8972293-0ef2-427e-b81c-380ee44cbdd1.md", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/417d83acef18db16e4504269791455a1b7e15ec5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.zh-cn.xlf", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c9255525979e9db716f4edd38e1c5885a9e23a62/e2e/f8972293-0ef2-427e-b81c-380ee44cbdd1.md", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74025245821fd71bdd61467e2aa2585ffc1af607/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.zh-cn.xlf", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/.localization-config", "", "", ".localization-config") | Out-Null

# ---- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.md"
$wsDe.Range("C2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.de-de.xlf"
$wsDe.Range("E2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.md"
$wsDe.Range("F2").Value = "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.de-de.xlf"

$wsDe.Range("A3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md"
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.de-de.xlf"
$wsDe.Range("D3").Value = "2016-01-26 11:53:30"
$wsDe.Range("E3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md"
$wsDe.Range("F3").Value = "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/e2e/937e6c38-aff9-4180-b7f0-59051d5eb1b8.md", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0bed570df1e0847441bf1ec621b6c2f8939eb0e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.de-de.xlf", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fdb065558eb5b15619c272213fb711b080ce2fc0/e2e/937e6c38-aff9-4180-b7f0-59051d5eb1b8.md", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/570cd94eab0c04779c06cde55542074602ae69ab/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.de-de.xlf", "", "", "f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/e2e/f8972293-0ef2-427e-b81c-380ee44cbdd1.md", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0bed570df1e0847441bf1ec621b6c2f8939eb0e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.de-de.xlf", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fdb065558eb5b15619c272213fb711b080ce2fc0/e2e/f8972293-0ef2-427e-b81c-380ee44cbdd1.md", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/570cd94eab0c04779c06cde55542074602ae69ab/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f8972293-0ef2-427e-b81c-380ee44cbdd1.26729d01a5c97ec124fed5697a2a96d96d78b6cf.de-de.xlf", "", "", "937e6c38-aff9-4180-b7f0-59051d5eb1b8.9128b6505a5fe5c4b64df1e5784ed96f4ad04442.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3ae4ff498e02d36d701eabd9dbd5cd0afbfcc9d9/.localization-config", "", "", ".localization-config") | Out-Null

$wb.Save()
